# Auto-generated edit script applying numeric corrections to the
# currentAveragePrice / LevePrice / LeveProfit columns (H:N) across
# several sheets, per the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 2150.6667
$ws.Range("I103").Value = 1976
$ws.Range("J103").Value = 2500
$ws.Range("K103").Value = 5928
$ws.Range("L103").Value = 7500
$ws.Range("M103").Value = -5342
$ws.Range("N103").Value = -8672

$ws.Range("H106").Value = 7160.577
$ws.Range("I106").Value = 2106.2307
$ws.Range("J106").Value = 12214.923
$ws.Range("K106").Value = 2106.2307
$ws.Range("L106").Value = 12214.923
$ws.Range("M106").Value = -1475.2307
$ws.Range("N106").Value = -13476.923

$ws.Range("H116").Value = 4404.5864
$ws.Range("I116").Value = 3931.111
$ws.Range("J116").Value = 5179.364
$ws.Range("K116").Value = 3931.111
$ws.Range("L116").Value = 5179.364
$ws.Range("M116").Value = -489.1109999999999
$ws.Range("N116").Value = -12063.364

$ws.Range("H121").Value = 1516796
$ws.Range("I121").Value = 3547.5
$ws.Range("J121").Value = 1684934.8
$ws.Range("K121").Value = 10642.5
$ws.Range("L121").Value = 5054804.4
$ws.Range("M121").Value = -8895.5
$ws.Range("N121").Value = -5058298.4

$ws.Range("H129").Value = 1390.6
$ws.Range("I129").Value = 915.1429000000001
$ws.Range("J129").Value = 2500
$ws.Range("K129").Value = 2745.4287
$ws.Range("L129").Value = 7500
$ws.Range("M129").Value = 2254.5713
$ws.Range("N129").Value = -17500

$ws.Range("H132").Value = 2991.9333
$ws.Range("I132").Value = 2990.923
$ws.Range("K132").Value = 8972.769
$ws.Range("M132").Value = -6442.769

$ws.Range("H137").Value = 12664.083
$ws.Range("I137").Value = 2000.625
$ws.Range("J137").Value = 17995.812
$ws.Range("K137").Value = 6001.875
$ws.Range("L137").Value = 53987.436
$ws.Range("M137").Value = -3451.875
$ws.Range("N137").Value = -59087.436

$ws.Range("H138").Value = 11754.263
$ws.Range("I138").Value = 9938.6
$ws.Range("J138").Value = 12402.714
$ws.Range("K138").Value = 29815.8
$ws.Range("L138").Value = 37208.142
$ws.Range("M138").Value = -24675.8
$ws.Range("N138").Value = -47488.142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5562941.5
$ws.Range("I32").Value = 5957286.5
$ws.Range("K32").Value = 5957286.5
$ws.Range("M32").Value = -5956999.5

$ws.Range("H45").Value = 2360.625
$ws.Range("I45").Value = 1966.3334
$ws.Range("J45").Value = 2451.6155
$ws.Range("K45").Value = 1966.3334
$ws.Range("L45").Value = 2451.6155
$ws.Range("M45").Value = -1589.3334
$ws.Range("N45").Value = -3205.6155

$ws.Range("H102").Value = 1693.8572
$ws.Range("I102").Value = 1825.5
$ws.Range("J102").Value = 904
$ws.Range("K102").Value = 1825.5
$ws.Range("L102").Value = 904
$ws.Range("M102").Value = -203.5
$ws.Range("N102").Value = -4148

$ws.Range("H132").Value = 5636.758
$ws.Range("I132").Value = 1512.96
$ws.Range("K132").Value = 4538.88
$ws.Range("M132").Value = -2008.88

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2078.8635
$ws.Range("I20").Value = 1972.3334
$ws.Range("J20").Value = 2558.25
$ws.Range("K20").Value = 1972.3334
$ws.Range("L20").Value = 2558.25
$ws.Range("M20").Value = -1725.3334
$ws.Range("N20").Value = -3052.25

$ws.Range("H134").Value = 39578.484
$ws.Range("I134").Value = 4777.15
$ws.Range("J134").Value = 116914.78
$ws.Range("K134").Value = 14331.45
$ws.Range("L134").Value = 350744.34
$ws.Range("M134").Value = -11796.45
$ws.Range("N134").Value = -355814.34

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 549.9231
$ws.Range("I19").Value = 490
$ws.Range("J19").Value = 749.6667
$ws.Range("K19").Value = 490
$ws.Range("L19").Value = 749.6667
$ws.Range("M19").Value = -320
$ws.Range("N19").Value = -1089.6667

$ws.Range("H24").Value = 549.9231
$ws.Range("I24").Value = 490
$ws.Range("J24").Value = 749.6667
$ws.Range("K24").Value = 490
$ws.Range("L24").Value = 749.6667
$ws.Range("M24").Value = -320
$ws.Range("N24").Value = -1089.6667

$ws.Range("H31").Value = 502796.53
$ws.Range("I31").Value = 16503.092
$ws.Range("J31").Value = 651386.2
$ws.Range("K31").Value = 16503.092
$ws.Range("L31").Value = 651386.2
$ws.Range("M31").Value = -16208.092
$ws.Range("N31").Value = -651976.2

$ws.Range("H34").Value = 502796.53
$ws.Range("I34").Value = 16503.092
$ws.Range("J34").Value = 651386.2
$ws.Range("K34").Value = 16503.092
$ws.Range("L34").Value = 651386.2
$ws.Range("M34").Value = -16301.092
$ws.Range("N34").Value = -651790.2

$ws.Range("H62").Value = 2356.9092
$ws.Range("I62").Value = 2115.875
$ws.Range("J62").Value = 2999.6667
$ws.Range("K62").Value = 2115.875
$ws.Range("L62").Value = 2999.6667
$ws.Range("M62").Value = -1491.875
$ws.Range("N62").Value = -4247.6667

$ws.Range("H65").Value = 2356.9092
$ws.Range("I65").Value = 2115.875
$ws.Range("J65").Value = 2999.6667
$ws.Range("K65").Value = 10579.375
$ws.Range("L65").Value = 14998.3335
$ws.Range("M65").Value = -7459.375
$ws.Range("N65").Value = -21238.3335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 555.75
$ws.Range("I14").Value = 555.75
$ws.Range("K14").Value = 1667.25
$ws.Range("M14").Value = -1494.25

$ws.Range("H31").Value = 5000
$ws.Range("I31").Value = 5000
$ws.Range("K31").Value = 15000
$ws.Range("M31").Value = -14712

$ws.Range("I33").Value = 883.3333
$ws.Range("J33").Value = 350
$ws.Range("K33").Value = 5299.9998
$ws.Range("L33").Value = 2100
$ws.Range("M33").Value = -5016.9998
$ws.Range("N33").Value = -2666

$ws.Range("H68").Value = 4271.4375
$ws.Range("I68").Value = 4379.5713
$ws.Range("J68").Value = 4241.16
$ws.Range("K68").Value = 13138.7139
$ws.Range("L68").Value = 12723.48
$ws.Range("M68").Value = -12327.7139
$ws.Range("N68").Value = -14345.48

$ws.Range("H71").Value = 4271.4375
$ws.Range("I71").Value = 4379.5713
$ws.Range("J71").Value = 4241.16
$ws.Range("K71").Value = 39416.14169999999
$ws.Range("L71").Value = 38170.44
$ws.Range("M71").Value = -35360.14169999999
$ws.Range("N71").Value = -46282.44

$ws.Range("H137").Value = 4849.5
$ws.Range("I137").Value = 4849.5
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 14548.5
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -9448.5
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6154
$ws.Range("I70").Value = 4231
$ws.Range("J70").Value = 10000
$ws.Range("K70").Value = 4231
$ws.Range("L70").Value = 10000
$ws.Range("M70").Value = -3961
$ws.Range("N70").Value = -10540

$ws.Range("H73").Value = 6154
$ws.Range("I73").Value = 4231
$ws.Range("J73").Value = 10000
$ws.Range("K73").Value = 4231
$ws.Range("L73").Value = 10000
$ws.Range("M73").Value = -3295
$ws.Range("N73").Value = -11872

$ws.Range("H97").Value = 1423.6522
$ws.Range("I97").Value = 1192.2
$ws.Range("J97").Value = 2966.6667
$ws.Range("K97").Value = 1192.2
$ws.Range("L97").Value = 2966.6667
$ws.Range("M97").Value = -696.2
$ws.Range("N97").Value = -3958.6667

$ws.Range("H102").Value = 4382.1113
$ws.Range("I102").Value = 3250.353
$ws.Range("J102").Value = 6306.1
$ws.Range("K102").Value = 3250.353
$ws.Range("L102").Value = 6306.1
$ws.Range("M102").Value = -1628.353
$ws.Range("N102").Value = -9550.1

$ws.Range("H122").Value = 7797.121
$ws.Range("I122").Value = 10079.15
$ws.Range("J122").Value = 4286.3076
$ws.Range("K122").Value = 30237.45
$ws.Range("L122").Value = 12858.9228
$ws.Range("M122").Value = -27787.45
$ws.Range("N122").Value = -17758.9228

$ws.Range("H132").Value = 38467450
$ws.Range("I132").Value = 45456576
$ws.Range("J132").Value = 27248.5
$ws.Range("K132").Value = 136369728
$ws.Range("L132").Value = 81745.5
$ws.Range("M132").Value = -136367198
$ws.Range("N132").Value = -86805.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1299
$ws.Range("I68").Value = 1299
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1299
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -550
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 1299
$ws.Range("I71").Value = 1299
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 6495
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -2751
$ws.Range("N71").ClearContents()

$ws.Range("H82").Value = 1285.95
$ws.Range("I82").Value = 855.7857
$ws.Range("J82").Value = 2289.6667
$ws.Range("K82").Value = 855.7857
$ws.Range("L82").Value = 2289.6667
$ws.Range("M82").Value = -494.7857
$ws.Range("N82").Value = -3011.6667

$ws.Range("H85").Value = 1285.95
$ws.Range("I85").Value = 855.7857
$ws.Range("J85").Value = 2289.6667
$ws.Range("K85").Value = 855.7857
$ws.Range("L85").Value = 2289.6667
$ws.Range("M85").Value = 392.2143
$ws.Range("N85").Value = -4785.6667

$ws.Range("H122").Value = 4144.2163
$ws.Range("I122").Value = 3385.5
$ws.Range("J122").Value = 9000
$ws.Range("K122").Value = 10156.5
$ws.Range("L122").Value = 27000
$ws.Range("M122").Value = -7706.5
$ws.Range("N122").Value = -31900

$ws.Range("H136").Value = 32052.158
$ws.Range("I136").Value = 6039.4614
$ws.Range("J136").Value = 69626.05499999999
$ws.Range("K136").Value = 18118.3842
$ws.Range("L136").Value = 208878.165
$ws.Range("M136").Value = -15568.3842
$ws.Range("N136").Value = -213978.165

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 15391150
$ws.Range("I62").Value = 5482
$ws.Range("K62").Value = 5482
$ws.Range("M62").Value = -4858

$ws.Range("H65").Value = 15391150
$ws.Range("I65").Value = 5482
$ws.Range("K65").Value = 27410
$ws.Range("M65").Value = -24290

$ws.Range("H107").Value = 26317756
$ws.Range("I107").Value = 33335708
$ws.Range("J107").Value = 438.25
$ws.Range("K107").Value = 100007124
$ws.Range("L107").Value = 1314.75
$ws.Range("M107").Value = -100005204
$ws.Range("N107").Value = -5154.75

$ws.Range("H132").Value = 3109.697
$ws.Range("I132").Value = 2549.1482
$ws.Range("J132").Value = 5632.1665
$ws.Range("K132").Value = 7647.444600000001
$ws.Range("L132").Value = 16896.4995
$ws.Range("M132").Value = -5117.444600000001
$ws.Range("N132").Value = -21956.4995

$ws.Range("H136").Value = 2022.875
$ws.Range("I136").Value = 2247.1667
$ws.Range("J136").Value = 1350
$ws.Range("K136").Value = 6741.500100000001
$ws.Range("L136").Value = 1350
$ws.Range("M136").Value = -4191.500100000001
$ws.Range("N136").Value = -9150

